$wb = $excel.ActiveWorkbook

# --- Palas sheet (sheet2): relabel the first 4 rows as PAA..PAD (regression
# buckets), drop the now-unneeded extra rows 6-11, and keep the shared
# formula for column I covering only the remaining rows.
$wsPalas = $wb.Worksheets.Item("Palas")
$wsPalas.Range("A2").Value = "PAA"
$wsPalas.Range("A3").Value = "PAB"
$wsPalas.Range("A4").Value = "PAC"
$wsPalas.Range("A5").Value = "PAD"
$wsPalas.Range("A6:J11").Delete()
$wsPalas.Range("I3:I5").Formula = "=70/60"

# --- Sheet1 (sheet3): add a header row on top of the existing data table.
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("A1").Value = "ID"
$wsSheet1.Range("B1").Value = "Media Caida (Exponencial) / min"
$wsSheet1.Range("C1").Value = "Media Recuperacion (Exponencial) / min"
$wsSheet1.Range("D1").Value = "NP"
$wsSheet1.Range("E1").Value = "NF+NP"
$wsSheet1.Range("F1").Value = "Min Carga"
$wsSheet1.Range("G1").Value = "Moda Carga"
$wsSheet1.Range("H1").Value = "Max Carga"
$wsSheet1.Range("I1").Value = "Aculatamiento"
$wsSheet1.Range("J1").Value = "Capacidad"

# --- View state: Palas is no longer the focused tab/selection, Sheet1 is.
$wsPalas.Activate()
$wsPalas.Range("A1:XFD1").Select()

$wsSheet1.Activate()
$wsSheet1.Range("F1:H1048576").Select()
